$p = $ppt.ActivePresentation

# Slide 4: "Sublime from the command line"
$s = $p.Slides.Item(4)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# Paragraph 3 was "Windows: " -> "Windows on the “cygwin Terminal”:"
$para3 = $tr.Paragraphs(3, 1)
$rep3 = $para3.Characters(8, 2)
$rep3.Text = " on the “cygwin Terminal”:"

# Paragraph 7 (Mac "ln -s ..." command) fix the subl symlink target for mac:
#   ...SharedSupport/bin/subl" ~/bin/subl  ->  ...SharedSupport/bin/subl” /usr/local/bin/subl
$para7 = $tr.Paragraphs(7, 1)
$quoteChar = $para7.Characters(72, 1)
$quoteChar.Text = "”"
$target = $para7.Characters(74, 1)
$target.Text = "/usr/local"
